$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Donovan Mitchell"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Cleveland Cavaliers"

# Row 4
$ws.Range("A4").Value = "Malik Beasley"
$ws.Range("B4").Value = "SG"
$ws.Range("C4").Value = "Detroit Pistons"

# Row 5
$ws.Range("A5").Value = "Dyson Daniels"
$ws.Range("B5").Value = "PG,SG"
$ws.Range("C5").Value = "Atlanta Hawks"

# Row 7
$ws.Range("A7").Value = "De'Andre Hunter"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Atlanta Hawks"

# Row 8
$ws.Range("A8").Value = "Toumani Camara"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Portland Trail Blazers"

# Row 9
$ws.Range("A9").Value = "Michael Porter Jr."
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Denver Nuggets"

# Row 10
$ws.Range("A10").Value = "Myles Turner"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Indiana Pacers"

# Row 11
$ws.Range("A11").Value = "Jonas Valanciunas"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Washington Wizards"

# Row 12
$ws.Range("A12").Value = "Domantas Sabonis"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Sacramento Kings"

# Row 13
$ws.Range("A13").Value = "Luguentz Dort"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Oklahoma City Thunder"

# Row 14
$ws.Range("A14").Value = "Alexandre Sarr"
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Washington Wizards"

# Row 15
$ws.Range("A15").Value = "Victor Wembanyama"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "San Antonio Spurs"

# Row 16
$ws.Range("A16").Value = "Tari Eason"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Houston Rockets"

# Row 19
$ws.Range("A19").Value = "Deandre Ayton"
$ws.Range("B19").Value = "C"
$ws.Range("C19").Value = "Portland Trail Blazers"
